# Update "想去人数" (number of people interested) figures for several
# events on the "展览" sheet and the merged "全部类型" sheet, matching
# the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 5026   # was 5022
$wsExhibit.Range("F7").Value  = 397    # was 396
$wsExhibit.Range("F8").Value  = 576    # was 575
$wsExhibit.Range("F13").Value = 544    # was 543
$wsExhibit.Range("F18").Value = 730    # was 729
$wsExhibit.Range("F27").Value = 1983   # was 1978
$wsExhibit.Range("F37").Value = 571    # was 568

# Sheet "全部类型" (All types) - same events, repeated rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value  = 5026   # was 5022
$wsAll.Range("F9").Value  = 397    # was 396
$wsAll.Range("F10").Value = 576    # was 575
$wsAll.Range("F19").Value = 544    # was 543
$wsAll.Range("F25").Value = 730    # was 729
$wsAll.Range("F34").Value = 1983   # was 1978
$wsAll.Range("F43").Value = 571    # was 568
